# Fied import investor access
# Rename the "CF1"/"CF2" custom-field header columns (S1/T1) to their
# proper labels, and leave the selection on S2 (matching the saved
# workbook state after the edit).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("T1").Value = "Custom field   2"
$ws.Range("S1").Value = "Custom Field 1"

$ws.Range("S2").Select() | Out-Null
